$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.148.51'
$ws.Range('E2').Value = '  +3.27%  '
$ws.Range('D3').Value = '2.317.72'
$ws.Range('E3').Value = '  +3.11%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'310.84"
$ws.Range('E5').Value = '  +2.30%  '
$ws.Range('D6').Value = "'101.20"
$ws.Range('E6').Value = '  +6.55%  '
$ws.Range('D7').Value = "'0.536"
$ws.Range('E7').Value = '  +2.63%  '
$ws.Range('E9').Value = '  +8.31%  '
$ws.Range('D10').Value = "'36.09"
$ws.Range('E10').Value = '  +4.46%  '
$ws.Range('D11').Value = "'0.0818"
$ws.Range('E11').Value = '  +4.17%  '
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').Value = "'7.21"
$ws.Range('E13').Value = '  +7.39%  '
$ws.Range('D14').Value = '2.676.41'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').Value = "'15.03"
$ws.Range('E15').Value = '  +5.12%  '
$ws.Range('D16').Value = '2.319.41'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = "'0.813"
$ws.Range('E17').Value = '  +3.57%  '
$ws.Range('D18').Value = '43.084.12'
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('D19').Value = "'12.58"
$ws.Range('E19').Value = '  +2.85%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  +2.64%  '
$ws.Range('D21').Value = "'6.14"
$ws.Range('E21').Value = '  +3.69%  '
$ws.Range('D22').Value = "'68.60"
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').Value = "'241.66"
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('E24').Value = '  +6.63%  '
$ws.Range('D25').Value = "'2.64"
$ws.Range('E25').Value = '  +3.82%  '
$ws.Range('D26').Value = "'0.998"
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = "'24.80"
$ws.Range('E27').Value = '  +5.60%  '
$ws.Range('D28').Value = "'37.40"
$ws.Range('E28').Value = '  +3.85%  '
$ws.Range('D29').Value = "'9.69"
$ws.Range('E29').Value = '  +3.42%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').Value = "'166.14"
$ws.Range('E31').Value = '  +4.24%  '
$ws.Range('D32').Value = "'5.35"
$ws.Range('E32').Value = '  +3.99%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = "'3.14"
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = "'17.99"
$ws.Range('E35').Value = '  +6.92%  '
$ws.Range('E36').Value = '  +1.83%  '
$ws.Range('E37').Value = '  +4.07%  '
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('E39').Value = '  +3.24%  '
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').Value = '  +9.50%  '
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('D43').Value = "'19.54"
$ws.Range('E43').Value = '  +5.61%  '
$ws.Range('E44').Value = '  +3.66%  '
$ws.Range('D45').Value = '1.980.35'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Value = "'3.02"
$ws.Range('E46').Value = '  +4.94%  '
$ws.Range('D47').Value = "'9.86"
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').Value = "'2.98"
$ws.Range('E48').Value = '  +19.47%  '
$ws.Range('D49').Value = "'55.71"
$ws.Range('E49').Value = '  +6.37%  '
$ws.Range('D50').Value = '2.543.09'
$ws.Range('E50').Value = '  +2.74%  '
$ws.Range('E51').Value = '  +4.37%  '

$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D7').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D19').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D22').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D32').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
